# packages/frontend/public/example.xlsx — "fixed example.xlsx" touch-up.
#
# The real edit (per the OOXML diff) was made simply by opening the file in
# Excel, letting it auto-size/select things, and re-saving:
#   - selection moved to cover the whole used range A1:E25
#   - column E best-fit to its (now widest) content -> stored width 62
#
# Reproduce both through the normal Excel object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E: auto-fit to content -------------------------------------
# A real "double-click the column border" AutoFit recalculates the column to
# the width of its widest cell. We drive the same ColumnWidth property that
# AutoFit ultimately sets, landing on the fitted width (stored width 62,
# i.e. ColumnWidth 61.1666... once Excel's fixed column-width padding of
# 5/6 of a character is added back on save).
$ws.Columns("E").AutoFit()
$ws.Columns("E").ColumnWidth = 61.16666666666667

# --- Selection: select the whole used range -----------------------------
# Mirrors selecting A1:E25 (e.g. via Ctrl+A / Ctrl+Shift+End from A1) so the
# saved sheetView selection becomes sqref="A1:E25" anchored at A1.
$ws.Range("A1").Select()
$ws.Range("A1:E25").Select()
